$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H33").Value = 848.8
$ws.Range("I33").Value = 844.44446
$ws.Range("K33").Value = 844.44446
$ws.Range("M33").Value = -615.44446
$ws.Range("H39").Value = 85.2
$ws.Range("I39").Value = 58.666668
$ws.Range("J39").Value = 125
$ws.Range("K39").Value = 176.000004
$ws.Range("L39").Value = 375
$ws.Range("M39").Value = 119.999996
$ws.Range("N39").Value = -967
$ws.Range("H42").Value = 11.6
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H46").Value = 3509
$ws.Range("J46").Value = 3509
$ws.Range("L46").Value = 10527
$ws.Range("N46").Value = -10765
$ws.Range("H51").Value = 9779.450000000001
$ws.Range("J51").Value = 10724.375
$ws.Range("L51").Value = 10724.375
$ws.Range("N51").Value = -11692.375
$ws.Range("H60").Value = 3509
$ws.Range("J60").Value = 3509
$ws.Range("L60").Value = 10527
$ws.Range("N60").Value = -11495
$ws.Range("H70").Value = 7109.222
$ws.Range("I70").Value = 10000
$ws.Range("J70").Value = 6747.875
$ws.Range("K70").Value = 30000
$ws.Range("L70").Value = 20243.625
$ws.Range("M70").Value = -29730
$ws.Range("N70").Value = -20783.625
$ws.Range("H73").Value = 7109.222
$ws.Range("I73").Value = 10000
$ws.Range("J73").Value = 6747.875
$ws.Range("K73").Value = 30000
$ws.Range("L73").Value = 20243.625
$ws.Range("M73").Value = -29064
$ws.Range("N73").Value = -22115.625
$ws.Range("H74").Value = 10212702
$ws.Range("I74").Value = 14291803
$ws.Range("K74").Value = 14291803
$ws.Range("M74").Value = -14290867
$ws.Range("H76").Value = 66671092
$ws.Range("I76").Value = 166669740
$ws.Range("J76").Value = 5318.3335
$ws.Range("K76").Value = 166669740
$ws.Range("L76").Value = 5318.3335
$ws.Range("M76").Value = -166669425
$ws.Range("N76").Value = -5948.3335
$ws.Range("H77").Value = 10212702
$ws.Range("I77").Value = 14291803
$ws.Range("K77").Value = 71459015
$ws.Range("M77").Value = -71454335
$ws.Range("H79").Value = 66671092
$ws.Range("I79").Value = 166669740
$ws.Range("J79").Value = 5318.3335
$ws.Range("K79").Value = 166669740
$ws.Range("L79").Value = 5318.3335
$ws.Range("M79").Value = -166668648
$ws.Range("N79").Value = -7502.3335
$ws.Range("H87").Value = 71851
$ws.Range("J87").Value = 71851
$ws.Range("L87").Value = 71851
$ws.Range("N87").Value = -74347
$ws.Range("H90").Value = 71851
$ws.Range("J90").Value = 71851
$ws.Range("L90").Value = 215553
$ws.Range("N90").Value = -228033
$ws.Range("H100").Value = 6310.5557
$ws.Range("I100").Value = 1756.4286
$ws.Range("J100").Value = 22250
$ws.Range("K100").Value = 1756.4286
$ws.Range("L100").Value = 22250
$ws.Range("M100").Value = -1215.4286
$ws.Range("N100").Value = -23332
$ws.Range("H104").Value = 94.75
$ws.Range("I104").Value = 76.333336
$ws.Range("J104").Value = 150
$ws.Range("K104").Value = 229.000008
$ws.Range("L104").Value = 450
$ws.Range("M104").Value = 1517.999992
$ws.Range("N104").Value = -3944
$ws.Range("H107").Value = 5201.6665
$ws.Range("I107").Value = 5222.7856
$ws.Range("J107").Value = 4906
$ws.Range("K107").Value = 5222.7856
$ws.Range("L107").Value = 4906
$ws.Range("M107").Value = -3302.7856
$ws.Range("N107").Value = -8746
$ws.Range("H138").Value = 4362.983
$ws.Range("J138").Value = 5408.6177
$ws.Range("L138").Value = 16225.8531
$ws.Range("N138").Value = -26505.8531

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H61").Value = 5459.7095
$ws.Range("I61").Value = 6077.1577
$ws.Range("K61").Value = 6077.1577
$ws.Range("M61").Value = -5865.1577
$ws.Range("H63").Value = 5829.857
$ws.Range("J63").Value = 6036.3335
$ws.Range("L63").Value = 6036.3335
$ws.Range("N63").Value = -7408.3335
$ws.Range("H66").Value = 5829.857
$ws.Range("J66").Value = 6036.3335
$ws.Range("L66").Value = 30181.6675
$ws.Range("N66").Value = -37045.6675
$ws.Range("H102").Value = 1317.2
$ws.Range("I102").Value = 1317.2
$ws.Range("K102").Value = 1317.2
$ws.Range("M102").Value = 304.8
$ws.Range("H110").Value = 5200.1904
$ws.Range("I110").Value = 2256.25
$ws.Range("K110").Value = 2256.25
$ws.Range("M110").Value = -211.25
$ws.Range("H136").Value = 5459.7095
$ws.Range("I136").Value = 6077.1577
$ws.Range("K136").Value = 18231.4731
$ws.Range("M136").Value = -15681.4731

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H20").Value = 6804295
$ws.Range("I20").Value = 8929644
$ws.Range("J20").Value = 3179.4
$ws.Range("K20").Value = 8929644
$ws.Range("L20").Value = 3179.4
$ws.Range("M20").Value = -8929397
$ws.Range("N20").Value = -3673.4
$ws.Range("H82").Value = 24217.428
$ws.Range("J82").Value = 33332.5
$ws.Range("L82").Value = 33332.5
$ws.Range("N82").Value = -34098.5
$ws.Range("H85").Value = 24217.428
$ws.Range("J85").Value = 33332.5
$ws.Range("L85").Value = 33332.5
$ws.Range("N85").Value = -35984.5
$ws.Range("H99").Value = 9054.629999999999
$ws.Range("I99").Value = 8401.020500000001
$ws.Range("J99").Value = 10055.469
$ws.Range("K99").Value = 8401.020500000001
$ws.Range("L99").Value = 10055.469
$ws.Range("M99").Value = -6903.020500000001
$ws.Range("N99").Value = -13051.469
$ws.Range("H134").Value = 1207152.9
$ws.Range("I134").Value = 1514092.8
$ws.Range("K134").Value = 4542278.4
$ws.Range("M134").Value = -4539743.4

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H58").Value = 13673.786
$ws.Range("I58").Value = 11545.167
$ws.Range("J58").Value = 15270.25
$ws.Range("K58").Value = 11545.167
$ws.Range("L58").Value = 15270.25
$ws.Range("M58").Value = -11342.167
$ws.Range("N58").Value = -15676.25
$ws.Range("H62").Value = 6949.7856
$ws.Range("I62").Value = 6124.9165
$ws.Range("K62").Value = 6124.9165
$ws.Range("M62").Value = -5500.9165
$ws.Range("H65").Value = 6949.7856
$ws.Range("I65").Value = 6124.9165
$ws.Range("K65").Value = 30624.5825
$ws.Range("M65").Value = -27504.5825
$ws.Range("H136").Value = 13673.786
$ws.Range("I136").Value = 11545.167
$ws.Range("J136").Value = 15270.25
$ws.Range("K136").Value = 34635.501
$ws.Range("L136").Value = 45810.75
$ws.Range("M136").Value = -32085.501
$ws.Range("N136").Value = -50910.75

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H4").Value = 1526.2693
$ws.Range("I4").Value = 1346.5106
$ws.Range("J4").Value = 3216
$ws.Range("K4").Value = 4039.5318
$ws.Range("L4").Value = 9648
$ws.Range("M4").Value = -3927.5318
$ws.Range("N4").Value = -9872
$ws.Range("H37").Value = 122466
$ws.Range("J37").Value = 122466
$ws.Range("L37").Value = 367398
$ws.Range("N37").Value = -367622
$ws.Range("H104").Value = 1358.8
$ws.Range("I104").Value = 448.5
$ws.Range("J104").Value = 5000
$ws.Range("K104").Value = 1345.5
$ws.Range("L104").Value = 15000
$ws.Range("M104").Value = 1275.5
$ws.Range("N104").Value = -20242
$ws.Range("H121").Value = 131859.3
$ws.Range("I121").Value = 1937.2858
$ws.Range("K121").Value = 5811.857400000001
$ws.Range("M121").Value = -4501.857400000001
$ws.Range("H140").Value = 34093710
$ws.Range("I140").Value = 39476300
$ws.Range("K140").Value = 118428900
$ws.Range("M140").Value = -118423720

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H18").Value = 746.6667
$ws.Range("I18").Value = 746.6667
$ws.Range("K18").Value = 746.6667
$ws.Range("M18").Value = -453.6667
$ws.Range("H70").Value = 6470.3887
$ws.Range("I70").Value = 6270.6665
$ws.Range("K70").Value = 6270.6665
$ws.Range("M70").Value = -6000.6665
$ws.Range("H73").Value = 6470.3887
$ws.Range("I73").Value = 6270.6665
$ws.Range("K73").Value = 6270.6665
$ws.Range("M73").Value = -5334.6665
$ws.Range("H80").Value = 31256878
$ws.Range("I80").Value = 142860590
$ws.Range("J80").Value = 7836.48
$ws.Range("K80").Value = 142860590
$ws.Range("L80").Value = 7836.48
$ws.Range("M80").Value = -142859592
$ws.Range("N80").Value = -9832.48
$ws.Range("H83").Value = 31256878
$ws.Range("I83").Value = 142860590
$ws.Range("J83").Value = 7836.48
$ws.Range("K83").Value = 714302950
$ws.Range("L83").Value = 39182.39999999999
$ws.Range("M83").Value = -714297958
$ws.Range("N83").Value = -49166.39999999999
$ws.Range("H132").Value = 8109.4165
$ws.Range("I132").Value = 8322.421
$ws.Range("K132").Value = 24967.263
$ws.Range("M132").Value = -22437.263

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H22").Value = 721.6
$ws.Range("I22").Value = 688.7
$ws.Range("J22").Value = 787.4
$ws.Range("K22").Value = 688.7
$ws.Range("L22").Value = 787.4
$ws.Range("M22").Value = -393.7
$ws.Range("N22").Value = -1377.4
$ws.Range("H27").Value = 721.6
$ws.Range("I27").Value = 688.7
$ws.Range("J27").Value = 787.4
$ws.Range("K27").Value = 688.7
$ws.Range("L27").Value = 787.4
$ws.Range("M27").Value = -581.7
$ws.Range("N27").Value = -1001.4
$ws.Range("H68").Value = 2676.8333
$ws.Range("I68").Value = 2212.2
$ws.Range("K68").Value = 2212.2
$ws.Range("M68").Value = -1463.2
$ws.Range("H71").Value = 2676.8333
$ws.Range("I71").Value = 2212.2
$ws.Range("K71").Value = 11061
$ws.Range("M71").Value = -7317
$ws.Range("H136").Value = 71438380
$ws.Range("I136").Value = 13063.167
$ws.Range("K136").Value = 39189.501
$ws.Range("M136").Value = -36639.501
